$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "42÷5=8, 2"
$t.Cell(1, 2).Range.Text = "42÷9=4, 6"
$t.Cell(1, 3).Range.Text = "45÷6=7, 3"
$t.Cell(1, 4).Range.Text = "88÷4=22, 0"
$t.Cell(1, 5).Range.Text = "69÷8=8, 5"

$t.Cell(5, 1).Range.Text = "24÷7=3, 3"
$t.Cell(5, 2).Range.Text = "67÷9=7, 4"
$t.Cell(5, 3).Range.Text = "21÷8=2, 5"
$t.Cell(5, 4).Range.Text = "52÷3=17, 1"
$t.Cell(5, 5).Range.Text = "55÷7=7, 6"

$t.Cell(9, 1).Range.Text = "45÷5=9, 0"
$t.Cell(9, 2).Range.Text = "33÷7=4, 5"
$t.Cell(9, 3).Range.Text = "85÷4=21, 1"
$t.Cell(9, 4).Range.Text = "89÷6=14, 5"
$t.Cell(9, 5).Range.Text = "59÷2=29, 1"

$t.Cell(13, 1).Range.Text = "85÷8=10, 5"
$t.Cell(13, 2).Range.Text = "38÷5=7, 3"
$t.Cell(13, 3).Range.Text = "63÷8=7, 7"
$t.Cell(13, 4).Range.Text = "41÷9=4, 5"
$t.Cell(13, 5).Range.Text = "44÷2=22, 0"

$t.Cell(17, 1).Range.Text = "41÷3=13, 2"
$t.Cell(17, 2).Range.Text = "21÷6=3, 3"
$t.Cell(17, 3).Range.Text = "69÷3=23, 0"
$t.Cell(17, 4).Range.Text = "68÷3=22, 2"
$t.Cell(17, 5).Range.Text = "82÷8=10, 2"

